$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.219.67'
$ws.Range("E2").Value = '  +1.72%  '
$ws.Range("D3").Value = '2.055.19'
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.32'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("E6").Value = '  +2.33%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.67'
$ws.Range("E8").Value = '  +3.87%  '
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.07'
$ws.Range("E10").Value = '  +1.27%  '
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("E12").Value = '  +1.27%  '
$ws.Range("D13").Value = '2.357.93'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.34'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.771'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.16'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = '2.055.02'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").Value = '37.149.11'
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.31'
$ws.Range("E20").Value = '  +12.79%  '
$ws.Range("E21").Value = '  +1.77%  '
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '224.79'
$ws.Range("E23").Value = '  +1.50%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  +1.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.37'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.94'
$ws.Range("E27").Value = '  +1.90%  '
$ws.Range("E28").Value = '  +6.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.78'
$ws.Range("E29").Value = '  +1.05%  '
$ws.Range("E30").Value = '  -4.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.03'
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E32").Value = '  -0.61%  '
$ws.Range("E33").Value = '  +1.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0612'
$ws.Range("E34").Value = '  +1.56%  '
$ws.Range("E36").Value = '  +5.82%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.82'
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.26'
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("E40").Value = '  -1.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.61'
$ws.Range("E41").Value = '  +14.33%  '
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("D43").Value = '1.484.08'
$ws.Range("E43").Value = '  +1.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.56'
$ws.Range("E44").Value = '  +2.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.17'
$ws.Range("E45").Value = '  +4.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0927'
$ws.Range("E46").Value = '  -1.89%  '
$ws.Range("E47").Value = '  +2.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.32'
$ws.Range("E48").Value = '  -2.68%  '
$ws.Range("E49").Value = '  +1.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.17'
$ws.Range("E50").Value = '  +4.25%  '
$ws.Range("E51").Value = '  +2.08%  '
